$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.977.37'
$ws.Range("E2").Value = '  -0.21%  '

$ws.Range("D3").Value = '2.304.75'
$ws.Range("E3").Value = '  +0.02%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.512'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.66%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.507'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.87'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0789'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.22'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.30%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.119'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.25%  '

$ws.Range("D15").Value = '2.664.41'
$ws.Range("E15").Value = '  +0.05%  '

$ws.Range("D16").Value = '2.308.26'
$ws.Range("E16").Value = '  -0.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.785'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.01%  '

$ws.Range("D18").Value = '42.914.97'
$ws.Range("E18").Value = '  -0.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.11%  '

$ws.Range("D20").Value = '0.0₃0904'
$ws.Range("E20").Value = '  -0.33%  '

$ws.Range("E21").Value = '  -1.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.03'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.63'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.46%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.50'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.06%  '

$ws.Range("E26").Value = '  +0.19%  '

$ws.Range("E27").Value = '  -0.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.92%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.74%  '

$ws.Range("E33").Value = '  +0.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.82'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.50%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.03'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.82%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.26'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.55%  '

$ws.Range("E37").Value = '  -1.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0691'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.26%  '

$ws.Range("E39").Value = '  -1.30%  '

$ws.Range("E40").Value = '  -1.19%  '

$ws.Range("E41").Value = '  -1.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.27%  '

$ws.Range("D43").Value = '2.009.64'
$ws.Range("E43").Value = '  -0.10%  '

$ws.Range("E44").Value = '  -1.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.01%  '

$ws.Range("E48").Value = '  -1.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.92'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.43%  '

$ws.Range("D51").Value = '2.531.65'
